$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 209; this shifts existing rows 209-217
# down to 210-218 (preserving their data/formatting unchanged).
$ws.Rows("209:209").Insert()

# Populate the newly inserted row 209 with the new weekly record.
$ws.Range("A209").Value2 = 10
$ws.Range("B209").Value2 = "Vega Modelo de Temuco"
$ws.Range("C209").Value2 = "La Araucanía"
$ws.Range("D209").Value2 = 44509
$ws.Range("E209").Value2 = 9
$ws.Range("F209").Value2 = 100112009
$ws.Range("G209").Value2 = "Acelga"
$ws.Range("H209").Value2 = "Sin especificar"
$ws.Range("I209").Value2 = "Primera"
$ws.Range("J209").Value2 = 30
$ws.Range("K209").Value2 = 8000
$ws.Range("L209").Value2 = 8000
$ws.Range("M209").Value2 = 8000
$ws.Range("N209").Value2 = "$/docena de atados (12 kilos)"
$ws.Range("O209").Value2 = "Provincia de Cautín"
$ws.Range("P209").Value2 = 667
$ws.Range("Q209").Value2 = 12
$ws.Range("R209").Value2 = "Hortaliza"
